# Atualizacao de bases das ligas, do dia: 31-03-2024 as 20:29
# Bosnia Herzegovina Premier Liga.xlsx
#
# 1) Two pairs of fixture rows had their detail columns (B..AC) swapped
#    while keeping the running id in column A fixed in place: 36<->37,
#    76<->77, 87<->88, 122<->123.
# 2) Five new fixture rows were appended at the end of the sheet (140..144).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap content of paired rows (B..AC), keep A (id/index) column fixed ---
# Rows 36 and 37
$ws.Range("B36").Value = 6864629
$ws.Range("C36").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D36").Value = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E36").Value = 45186.61458333334
$ws.Range("F36").Value = "Borac Banja Luka"
$ws.Range("G36").Value = "NK Posusje"
$ws.Range("H36").Value = 1
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = "H"
$ws.Range("K36").Value = 1.363
$ws.Range("L36").Value = 4.5
$ws.Range("M36").Value = 6.5
$ws.Range("N36").Value = 1.363
$ws.Range("O36").Value = 4.2
$ws.Range("P36").Value = 6.5
$ws.Range("Q36").Value = -1.25
$ws.Range("R36").Value = 1.95
$ws.Range("S36").Value = 1.85
$ws.Range("T36").Value = 2.5
$ws.Range("U36").Value = 1.925
$ws.Range("V36").Value = 1.875
$ws.Range("W36").Value = 0.363
$ws.Range("X36").Value = -1
$ws.Range("Y36").Value = -1
$ws.Range("Z36").Value = -0.5
$ws.Range("AA36").Value = 0.425
$ws.Range("AB36").Value = -1
$ws.Range("AC36").Value = 0.875

$ws.Range("B37").Value = 6865299
$ws.Range("C37").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D37").Value = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E37").Value = 45186.61458333334
$ws.Range("F37").Value = "Siroki Brijeg"
$ws.Range("G37").Value = "Zvijezda 09"
$ws.Range("H37").Value = 2
$ws.Range("I37").Value = 1
$ws.Range("J37").Value = "H"
$ws.Range("K37").Value = 1.25
$ws.Range("L37").Value = 5.5
$ws.Range("M37").Value = 8
$ws.Range("N37").Value = 1.4
$ws.Range("O37").Value = 4.75
$ws.Range("P37").Value = 5.75
$ws.Range("Q37").Value = -1.25
$ws.Range("R37").Value = 1.9
$ws.Range("S37").Value = 1.9
$ws.Range("T37").Value = 2.75
$ws.Range("U37").Value = 1.85
$ws.Range("V37").Value = 1.95
$ws.Range("W37").Value = 0.3999999999999999
$ws.Range("X37").Value = -1
$ws.Range("Y37").Value = -1
$ws.Range("Z37").Value = -0.5
$ws.Range("AA37").Value = 0.45
$ws.Range("AB37").Value = 0.425
$ws.Range("AC37").Value = -0.5

# Rows 76 and 77
$ws.Range("B76").Value = 6865328
$ws.Range("C76").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D76").Value = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E76").Value = 45235.375
$ws.Range("F76").Value = "Siroki Brijeg"
$ws.Range("G76").Value = "NK Posusje"
$ws.Range("H76").Value = 1
$ws.Range("I76").Value = 1
$ws.Range("J76").Value = "D"
$ws.Range("K76").Value = 2
$ws.Range("L76").Value = 3
$ws.Range("M76").Value = 3.5
$ws.Range("N76").Value = 2.1
$ws.Range("O76").Value = 3
$ws.Range("P76").Value = 3.3
$ws.Range("Q76").Value = -0.25
$ws.Range("R76").Value = 1.825
$ws.Range("S76").Value = 1.975
$ws.Range("T76").Value = 2
$ws.Range("U76").Value = 1.825
$ws.Range("V76").Value = 1.975
$ws.Range("W76").Value = -1
$ws.Range("X76").Value = 2
$ws.Range("Y76").Value = -1
$ws.Range("Z76").Value = -0.5
$ws.Range("AA76").Value = 0.4875
$ws.Range("AB76").Value = 0
$ws.Range("AC76").Value = -0

$ws.Range("B77").Value = 6865377
$ws.Range("C77").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D77").Value = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E77").Value = 45235.375
$ws.Range("F77").Value = "Zrinjski Mostar"
$ws.Range("G77").Value = "FK Tuzla City"
$ws.Range("H77").Value = 3
$ws.Range("I77").Value = 1
$ws.Range("J77").Value = "H"
$ws.Range("K77").Value = 1.333
$ws.Range("L77").Value = 5
$ws.Range("M77").Value = 6
$ws.Range("N77").Value = 1.166
$ws.Range("O77").Value = 6.5
$ws.Range("P77").Value = 13
$ws.Range("Q77").Value = -2
$ws.Range("R77").Value = 1.9
$ws.Range("S77").Value = 1.9
$ws.Range("T77").Value = 3.25
$ws.Range("U77").Value = 1.95
$ws.Range("V77").Value = 1.85
$ws.Range("W77").Value = 0.1659999999999999
$ws.Range("X77").Value = -1
$ws.Range("Y77").Value = -1
$ws.Range("Z77").Value = 0
$ws.Range("AA77").Value = -0
$ws.Range("AB77").Value = 0.95
$ws.Range("AC77").Value = -1

# Rows 87 and 88
$ws.Range("B87").Value = 7505495
$ws.Range("C87").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D87").Value = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E87").Value = 45256.375
$ws.Range("F87").Value = "Sloga"
$ws.Range("G87").Value = "Zvijezda 09"
$ws.Range("H87").Value = 1
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = "H"
$ws.Range("K87").Value = 1.444
$ws.Range("L87").Value = 4.2
$ws.Range("M87").Value = 5.5
$ws.Range("N87").Value = 1.5
$ws.Range("O87").Value = 4.2
$ws.Range("P87").Value = 5.25
$ws.Range("Q87").Value = -1
$ws.Range("R87").Value = 1.8
$ws.Range("S87").Value = 2
$ws.Range("T87").Value = 2.75
$ws.Range("U87").Value = 1.775
$ws.Range("V87").Value = 2.025
$ws.Range("W87").Value = 0.5
$ws.Range("X87").Value = -1
$ws.Range("Y87").Value = -1
$ws.Range("Z87").Value = 0
$ws.Range("AA87").Value = -0
$ws.Range("AB87").Value = -1
$ws.Range("AC87").Value = 1.025

$ws.Range("B88").Value = 7505497
$ws.Range("C88").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D88").Value = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E88").Value = 45256.375
$ws.Range("F88").Value = "Zeljeznicar"
$ws.Range("G88").Value = "NK Posusje"
$ws.Range("H88").Value = 1
$ws.Range("I88").Value = 1
$ws.Range("J88").Value = "D"
$ws.Range("K88").Value = 1.65
$ws.Range("L88").Value = 3.4
$ws.Range("M88").Value = 4.75
$ws.Range("N88").Value = 1.8
$ws.Range("O88").Value = 3.2
$ws.Range("P88").Value = 4.2
$ws.Range("Q88").Value = -0.5
$ws.Range("R88").Value = 1.825
$ws.Range("S88").Value = 1.975
$ws.Range("T88").Value = 2
$ws.Range("U88").Value = 1.75
$ws.Range("V88").Value = 2.05
$ws.Range("W88").Value = -1
$ws.Range("X88").Value = 2.2
$ws.Range("Y88").Value = -1
$ws.Range("Z88").Value = -1
$ws.Range("AA88").Value = 0.9750000000000001
$ws.Range("AB88").Value = 0
$ws.Range("AC88").Value = -0

# Rows 122 and 123
$ws.Range("B122").Value = 6865363
$ws.Range("C122").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D122").Value = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E122").Value = 45353.375
$ws.Range("F122").Value = "NK Igman Konjic"
$ws.Range("G122").Value = "Siroki Brijeg"
$ws.Range("H122").Value = 1
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = "H"
$ws.Range("K122").Value = 2
$ws.Range("L122").Value = 3.3
$ws.Range("M122").Value = 3.25
$ws.Range("N122").Value = 2.3
$ws.Range("O122").Value = 3.2
$ws.Range("P122").Value = 2.75
$ws.Range("Q122").Value = -0.25
$ws.Range("R122").Value = 2.05
$ws.Range("S122").Value = 1.75
$ws.Range("T122").Value = 2
$ws.Range("U122").Value = 1.9
$ws.Range("V122").Value = 1.9
$ws.Range("W122").Value = 1.3
$ws.Range("X122").Value = -1
$ws.Range("Y122").Value = -1
$ws.Range("Z122").Value = 1.05
$ws.Range("AA122").Value = -1
$ws.Range("AB122").Value = -1
$ws.Range("AC122").Value = 0.8999999999999999

$ws.Range("B123").Value = 6865381
$ws.Range("C123").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D123").Value = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E123").Value = 45353.375
$ws.Range("F123").Value = "FK Tuzla City"
$ws.Range("G123").Value = "Zvijezda 09"
$ws.Range("H123").Value = 2
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = "H"
$ws.Range("K123").Value = 1.666
$ws.Range("L123").Value = 3.6
$ws.Range("M123").Value = 4.333
$ws.Range("N123").Value = 1.5
$ws.Range("O123").Value = 4
$ws.Range("P123").Value = 5.25
$ws.Range("Q123").Value = -1
$ws.Range("R123").Value = 1.925
$ws.Range("S123").Value = 1.875
$ws.Range("T123").Value = 2.5
$ws.Range("U123").Value = 1.8
$ws.Range("V123").Value = 2
$ws.Range("W123").Value = 0.5
$ws.Range("X123").Value = -1
$ws.Range("Y123").Value = -1
$ws.Range("Z123").Value = 0.925
$ws.Range("AA123").Value = -1
$ws.Range("AB123").Value = -1
$ws.Range("AC123").Value = 1

# --- Append new rows 140..144 (ids 138..142) ---
# Row 140
$ws.Range("A140").Value = 138
$ws.Range("B140").Value = 7952736
$ws.Range("C140").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D140").Value = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E140").Value = 45380.41666666666
$ws.Range("F140").Value = "NK Igman Konjic"
$ws.Range("G140").Value = "Sloga"
$ws.Range("H140").Value = 2
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = "H"
$ws.Range("K140").Value = 1.727
$ws.Range("L140").Value = 4
$ws.Range("M140").Value = 3.5
$ws.Range("N140").Value = 1.909
$ws.Range("O140").Value = 4
$ws.Range("P140").Value = 2.9
$ws.Range("Q140").Value = -0.5
$ws.Range("R140").Value = 2
$ws.Range("S140").Value = 1.8
$ws.Range("T140").Value = 2.5
$ws.Range("U140").Value = 1.9
$ws.Range("V140").Value = 1.9
$ws.Range("W140").Value = 0.909
$ws.Range("X140").Value = -1
$ws.Range("Y140").Value = -1
$ws.Range("Z140").Value = 1
$ws.Range("AA140").Value = -1
$ws.Range("AB140").Value = -1
$ws.Range("AC140").Value = 0.8999999999999999

# Row 141
$ws.Range("A141").Value = 139
$ws.Range("B141").Value = 7952738
$ws.Range("C141").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D141").Value = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E141").Value = 45380.5
$ws.Range("F141").Value = "FK Tuzla City"
$ws.Range("G141").Value = "FK Sarajevo"
$ws.Range("H141").Value = 1
$ws.Range("I141").Value = 5
$ws.Range("J141").Value = "A"
$ws.Range("K141").Value = 2.9
$ws.Range("L141").Value = 4
$ws.Range("M141").Value = 1.909
$ws.Range("N141").Value = 2.6
$ws.Range("O141").Value = 3.5
$ws.Range("P141").Value = 2.375
$ws.Range("Q141").Value = 0
$ws.Range("R141").Value = 1.95
$ws.Range("S141").Value = 1.85
$ws.Range("T141").Value = 2.5
$ws.Range("U141").Value = 1.9
$ws.Range("V141").Value = 1.9
$ws.Range("W141").Value = -1
$ws.Range("X141").Value = -1
$ws.Range("Y141").Value = 1.375
$ws.Range("Z141").Value = -1
$ws.Range("AA141").Value = 0.8500000000000001
$ws.Range("AB141").Value = 0.8999999999999999
$ws.Range("AC141").Value = -1

# Row 142
$ws.Range("A142").Value = 140
$ws.Range("B142").Value = 7952735
$ws.Range("C142").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D142").Value = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E142").Value = 45380.66666666666
$ws.Range("F142").Value = "Zeljeznicar"
$ws.Range("G142").Value = "Velez Mostar"
$ws.Range("H142").Value = 0
$ws.Range("I142").Value = 1
$ws.Range("J142").Value = "A"
$ws.Range("K142").Value = 2.4
$ws.Range("L142").Value = 3.2
$ws.Range("M142").Value = 2.6
$ws.Range("N142").Value = 2.5
$ws.Range("O142").Value = 3.2
$ws.Range("P142").Value = 2.625
$ws.Range("Q142").Value = 0
$ws.Range("R142").Value = 1.775
$ws.Range("S142").Value = 2.025
$ws.Range("T142").Value = 1.75
$ws.Range("U142").Value = 1.825
$ws.Range("V142").Value = 1.975
$ws.Range("W142").Value = -1
$ws.Range("X142").Value = -1
$ws.Range("Y142").Value = 1.625
$ws.Range("Z142").Value = -1
$ws.Range("AA142").Value = 1.025
$ws.Range("AB142").Value = -1
$ws.Range("AC142").Value = 0.9750000000000001

# Row 143
$ws.Range("A143").Value = 141
$ws.Range("B143").Value = 7952739
$ws.Range("C143").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D143").Value = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E143").Value = 45381.39583333334
$ws.Range("F143").Value = "Zvijezda 09"
$ws.Range("G143").Value = "Siroki Brijeg"
$ws.Range("H143").Value = 2
$ws.Range("I143").Value = 3
$ws.Range("J143").Value = "A"
$ws.Range("K143").Value = 2.25
$ws.Range("L143").Value = 3.1
$ws.Range("M143").Value = 2.875
$ws.Range("N143").Value = 2.15
$ws.Range("O143").Value = 2.75
$ws.Range("P143").Value = 3.4
$ws.Range("Q143").Value = -0.25
$ws.Range("R143").Value = 1.9
$ws.Range("S143").Value = 1.9
$ws.Range("T143").Value = 1.75
$ws.Range("U143").Value = 1.8
$ws.Range("V143").Value = 2
$ws.Range("W143").Value = -1
$ws.Range("X143").Value = -1
$ws.Range("Y143").Value = 2.4
$ws.Range("Z143").Value = -1
$ws.Range("AA143").Value = 0.8999999999999999
$ws.Range("AB143").Value = 0.8
$ws.Range("AC143").Value = -1

# Row 144
$ws.Range("A144").Value = 142
$ws.Range("B144").Value = 7952456
$ws.Range("C144").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D144").Value = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E144").Value = 45381.5
$ws.Range("F144").Value = "Borac Banja Luka"
$ws.Range("G144").Value = "Zrinjski Mostar"
$ws.Range("H144").Value = 1
$ws.Range("I144").Value = 2
$ws.Range("J144").Value = "A"
$ws.Range("K144").Value = 2.2
$ws.Range("L144").Value = 3.2
$ws.Range("M144").Value = 2.875
$ws.Range("N144").Value = 2.45
$ws.Range("O144").Value = 2.8
$ws.Range("P144").Value = 2.875
$ws.Range("Q144").Value = 0
$ws.Range("R144").Value = 1.725
$ws.Range("S144").Value = 2.075
$ws.Range("T144").Value = 1.75
$ws.Range("U144").Value = 1.75
$ws.Range("V144").Value = 2.05
$ws.Range("W144").Value = -1
$ws.Range("X144").Value = -1
$ws.Range("Y144").Value = 1.875
$ws.Range("Z144").Value = -1
$ws.Range("AA144").Value = 1.075
$ws.Range("AB144").Value = 0.75
$ws.Range("AC144").Value = -1

# Apply the same formatting used by the other data rows:
# column A = bold + thin border + centered/top aligned (style used for the id column)
# column E = custom date/time number format
$rA = $ws.Range("A140:A144")
$rA.Borders.LineStyle = 1
$rA.Font.Bold = $true
$rA.HorizontalAlignment = -4108
$rA.VerticalAlignment = -4160

$rE = $ws.Range("E140:E144")
$rE.NumberFormat = "YYYY-MM-DD HH:MM:SS"

